# MW to MWh correction
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values (columns B:F, rows 2:9) ---
$data = @(
    @(1.08240472833333, 550.96574363693298, 416.16601062843699, 0.59636792606230693, 0.45046005767583996),
    @(388.038161020833, 22.6838772985903, 16.377724522323899, 8.8022100317672098, 6.3551821053483595),
    @(1551.00895077833, 10.9604839706953, 8.2663260448639608, 16.999808743410799, 12.821145685636099),
    @(71.982325360000004, 97.219670258001699, 49.959810688289799, 6.9980979359034006, 3.59622334788848),
    @(1032.84026755083, 81.954842242142107, 74.477724496221001, 84.646261188460301, 76.923592895254103),
    @(1.2871983091666701, 884.29563198907294, 439.70755741240498, 1.1382638422997999, 0.56599082442905302),
    @(56.044908320000005, 71.834835441691197, -93.892694186698094, 4.02597676651187, -5.2622074376112904),
    @(1760.0948655091702, 669.35195732655097, 29.255740162575801, 1178.12294330897, 51.492878046819904)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

# --- Add new totals row 10 with SUM formulas ---
$ws.Range("E10").Formula = "=SUM(E2:E9)"
$ws.Range("F10").Formula = "=SUM(F2:F9)"
$ws.Range("E10:F10").NumberFormat = "#,##0.00_ ;[Red]\-#,##0.00\ "

# --- Number formats ---
# Columns B:D now show integers (0 decimals)
$ws.Range("B2:D9").NumberFormat = "#,##0_ ;[Red]\-#,##0\ "
# Columns E:F keep 2 decimals
$ws.Range("E2:F9").NumberFormat = "#,##0.00_ ;[Red]\-#,##0.00\ "

# --- Header row: remove bold/underline formatting, revert to default style ---
$ws.Rows("1:1").ClearFormats()

# --- Column widths: let Excel auto-fit based on new content ---
$ws.Columns("A:F").AutoFit()

# --- Selection / active cell ---
$ws.Range("E9").Select()
